# Natmi following Dr Hou advice
#
# The LR-pair results table (Slit2 -> Gpc1) is recomputed to include the
# "ECs" sending cluster alongside "FAPs" and "sCs": rows 2-4 now report the
# FAPs/Slit2/Gpc1 edge re-scored against an ECs/FAPs/sCs target cluster,
# rows 5-7 keep the (reordered) FAPs sending-cluster block, and three new
# rows (8-10) are appended for the "sCs" sending cluster. All of the
# per-edge statistics (detection rates, average/total expression, derived
# specificities, …) are refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> full row of values (A..T), in the final, post-edit layout.
$rowData = @{
    2  = @("ECs",  "Slit2", "Gpc1", "ECs",  2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 1.564139666666667, 4.692419,           0.166125338305886,  0.166125338305886,  0.05229440285555556, 0.4706496257,        0.0029126083891251,   0.002912608389125099)
    3  = @("ECs",  "Slit2", "Gpc1", "FAPs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 3.316850333333333, 9.950551000000001,  0.3522785691569683, 0.3522785691569683, 0.1108933628111111,  0.9980402653000001,  0.006176357720616414, 0.006176357720616412)
    4  = @("ECs",  "Slit2", "Gpc1", "sCs",  2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 4.534428999999999, 13.603287,           0.4815960925371456, 0.4815960925371456, 0.1516010762333333,  1.3644096861,        0.008443629572695108, 0.008443629572695106)
    5  = @("FAPs", "Slit2", "Gpc1", "ECs",  3, 1,                  0.37963,             1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 1.564139666666667, 4.692419,           0.166125338305886,  0.166125338305886,  0.5937943416566667,  5.34414907491,       0.03307218911556017,  0.03307218911556017)
    6  = @("FAPs", "Slit2", "Gpc1", "FAPs", 3, 1,                  0.37963,             1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 3.316850333333333, 9.950551000000001,  0.3522785691569683, 0.3522785691569683, 1.259175892043333,   11.33258302839,      0.07013152586672808,  0.07013152586672808)
    7  = @("FAPs", "Slit2", "Gpc1", "sCs",  3, 1,                  0.37963,             1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 4.534428999999999, 13.603287,           0.4815960925371456, 0.4815960925371456, 1.72140528127,       15.49264753143,      0.0958760247661688,   0.0958760247661688)
    8  = @("sCs",  "Slit2", "Gpc1", "ECs",  3, 1,                  1.493861,            4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 1.564139666666667, 4.692419, 0.166125338305886,  0.166125338305886,  2.336607246586333,   21.029465219277,     0.1301405408012007,   0.1301405408012007)
    9  = @("sCs",  "Slit2", "Gpc1", "FAPs", 3, 1,                  1.493861,            4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 3.316850333333333, 9.950551000000001, 0.3522785691569683, 0.3522785691569683, 4.954913355803667, 44.59422020223301, 0.2759706855696238,   0.2759706855696238)
    10 = @("sCs",  "Slit2", "Gpc1", "sCs",  3, 1,                  1.493861,            4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 4.534428999999999, 13.603287, 0.4815960925371456, 0.4815960925371456, 6.773806640368999,   60.964259763321,     0.3772764381982818,   0.3772764381982818)
}

foreach ($r in ($rowData.Keys | Sort-Object)) {
    $vals = $rowData[$r]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
